$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value that was bumped by one day
# (2023-09-15 -> 2023-09-16, serial 45184 -> 45185) for every data row
# (rows 2 through 28).
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
